$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'247.68"
$ws.Range("D3").Value = "'21.98"
$ws.Range("D4").Value = "'5.380"
$ws.Range("D5").Value = "'0.05638"
$ws.Range("D6").Value = "'3.432"
$ws.Range("D7").Value = "'6.354"
$ws.Range("D8").Value = "'0.8177"
$ws.Range("D9").Value = "'0.9346"
$ws.Range("D10").Value = "'0.1433"
$ws.Range("D11").Value = "'0.07445"
$ws.Range("D12").Value = "'0.03255"
$ws.Range("D13").Value = "'0.03091"
$ws.Range("D14").Value = "'0.09324"
$ws.Range("D15").Value = "'3.558"
$ws.Range("D17").Value = "'0.04737"
$ws.Range("D18").Value = "'0.0005791"
$ws.Range("D19").Value = "'0.006408"
$ws.Range("D20").Value = "'0.005064"
$ws.Range("D23").Value = "'3.751"
$ws.Range("D25").Value = "'0.3306"
$ws.Range("D26").Value = "'0.1319"
$ws.Range("D28").Value = "'0.0003000"
$ws.Range("D40").Value = "'0.03940"
$ws.Range("D41").Value = "'0.006864"
$ws.Range("D42").Value = "'0.1066"
$ws.Range("D43").Value = "'0.003022"
$ws.Range("D44").Value = "'0.008521"
$ws.Range("D45").Value = "'0.00005573"
$ws.Range("D48").Value = "'0.7801"
$ws.Range("D49").Value = "'0.1799"
$ws.Range("D50").Value = "'0.00002100"
